$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country re-ranking: Marruecos rises above Paises Bajos and Catar ---
$ws.Range("A33").Value = "Marruecos"
$ws.Range("A34").Value = "Paises Bajos"
$ws.Range("A35").Value = "Catar"

# --- Updated case statistics (B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados, G=Muertes hoy, H=Muertes) ---
$ws.Range("B4").Value = 7531397
$ws.Range("C4").Value = 36726
$ws.Range("D4").Value = 4762450
$ws.Range("E4").Value = 2555642
$ws.Range("G4").Value = 645
$ws.Range("H4").Value = 213305
$ws.Range("B5").Value = 6471734
$ws.Range("C5").Value = 79774
$ws.Range("D5").Value = 5424943
$ws.Range("E5").Value = 945918
$ws.Range("G5").Value = 1069
$ws.Range("H5").Value = 100873
$ws.Range("B25").Value = 298362
$ws.Range("C25").Value = 2832
$ws.Range("E25").Value = 29266
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 9596
$ws.Range("B33").Value = 128565
$ws.Range("C33").Value = 2521
$ws.Range("D33").Value = 106044
$ws.Range("E33").Value = 20258
$ws.Range("G33").Value = 34
$ws.Range("H33").Value = 2263
$ws.Range("B34").Value = 127922
$ws.Range("C34").Value = 3825
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("G34").Value = 9
$ws.Range("H34").Value = 6428
$ws.Range("B35").Value = 126164
$ws.Range("C35").Value = 205
$ws.Range("D35").Value = 123108
$ws.Range("E35").Value = 2841
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 215
$ws.Range("D68").Value = 18739
$ws.Range("E68").Value = 23048
$ws.Range("B103").Value = 10700
$ws.Range("C103").Value = 48
$ws.Range("D103").Value = 10039
$ws.Range("E103").Value = 595
$ws.Range("B109").Value = 8792
$ws.Range("C109").Value = 11
$ws.Range("D109").Value = 6949
$ws.Range("E109").Value = 1614
$ws.Range("B117").Value = 6205
$ws.Range("C117").Value = 79
$ws.Range("D117").Value = 5399
$ws.Range("E117").Value = 744
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 62
$ws.Range("B120").Value = 5718
$ws.Range("C120").Value = 48
$ws.Range("E120").Value = 674
$ws.Range("B136").Value = 4038
$ws.Range("C136").Value = 40
$ws.Range("D136").Value = 3406
$ws.Range("E136").Value = 602
$ws.Range("G136").Value = 3
$ws.Range("H136").Value = 30
$ws.Range("B152").Value = 2252
$ws.Range("C152").Value = 14
$ws.Range("D152").Value = 1697
$ws.Range("E152").Value = 483
$ws.Range("B166").Value = 1211
$ws.Range("C166").Value = 8
$ws.Range("D166").Value = 1070
$ws.Range("E166").Value = 56
$ws.Range("B176").Value = 513
$ws.Range("C176").Value = 3
$ws.Range("E176").Value = 40
